$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 'T'
$ws.Range("E2").Value = 'Let''s look at another STUDENT 1nswer.'
$ws.Range("D3").Value = 'T'
$ws.Range("D4").Value = 'T'
$ws.Range("D5").Value = 'T'
$ws.Range("D6").Value = 'T'
$ws.Range("D7").Value = 'S'
$ws.Range("D8").Value = 'S'
$ws.Range("D9").Value = 'SN'
$ws.Range("D10").Value = 'S'
$ws.Range("D13").Value = 'SN'
$ws.Range("D15").Value = 'SN'
$ws.Range("D16").Value = 'SN'
$ws.Range("D17").Value = 'S'
$ws.Range("D18").Value = 'T'
$ws.Range("D19").Value = 'T'
$ws.Range("D21").Value = 'T'
$ws.Range("D22").Value = 'T'
$ws.Range("D24").Value = 'T'
$ws.Range("D25").Value = 'T'
$ws.Range("D28").Value = 'T'
$ws.Range("D29").Value = 'T'
$ws.Range("D31").Value = 'T'
$ws.Range("D32").Value = 'T'
$ws.Range("D35").Value = 'T'
$ws.Range("D36").Value = 'T'
$ws.Range("D37").Value = 'T'
$ws.Range("D38").Value = 'T'
$ws.Range("D39").Value = 'T'
$ws.Range("D43").Value = 'T'
$ws.Range("D47").Value = 'T'
$ws.Range("D48").Value = 'T'
$ws.Range("D50").Value = 'T'
$ws.Range("D51").Value = 'T'
$ws.Range("D52").Value = 'T'
$ws.Range("D54").Value = 'T'
$ws.Range("D56").Value = 'T'
$ws.Range("D57").Value = 'T'
$ws.Range("D58").Value = 'T'
$ws.Range("D60").Value = 'T'
$ws.Range("D62").Value = 'T'
$ws.Range("D64").Value = 'T'
$ws.Range("D66").Value = 'T'
$ws.Range("D67").Value = 'T'
$ws.Range("D69").Value = 'T'
$ws.Range("D73").Value = 'T'
$ws.Range("D74").Value = 'T'
$ws.Range("D75").Value = 'T'
$ws.Range("D81").Value = 'T'
$ws.Range("D85").Value = 'T'
$ws.Range("D86").Value = 'T'
$ws.Range("D87").Value = 'T'
$ws.Range("D89").Value = 'T'
$ws.Range("D91").Value = 'T'
$ws.Range("D93").Value = 'T'
$ws.Range("D96").Value = 'T'
$ws.Range("D97").Value = 'T'
$ws.Range("D99").Value = 'T'
$ws.Range("D100").Value = 'T'
$ws.Range("D102").Value = 'T'
$ws.Range("D103").Value = 'T'
$ws.Range("D104").Value = 'S'
$ws.Range("D105").Value = 'T'
